$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style used by the
# other header cells (B1:E1) -- copy formatting from E1.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Fill F2:F65 with the recorded time_taken timestamps for each data row.
$timestamps = @(
    "2021-10-05 13:39:27.593466",
    "2021-10-05 13:39:27.593480",
    "2021-10-05 13:39:27.593484",
    "2021-10-05 13:39:27.593487",
    "2021-10-05 13:39:27.593491",
    "2021-10-05 13:39:27.593494",
    "2021-10-05 13:39:27.593497",
    "2021-10-05 13:39:27.593500",
    "2021-10-05 13:39:27.593503",
    "2021-10-05 13:39:27.593506",
    "2021-10-05 13:39:27.593510",
    "2021-10-05 13:39:27.593513",
    "2021-10-05 13:39:27.593516",
    "2021-10-05 13:39:27.593519",
    "2021-10-05 13:39:27.593522",
    "2021-10-05 13:39:27.593525",
    "2021-10-05 13:39:27.593528",
    "2021-10-05 13:39:27.593531",
    "2021-10-05 13:39:27.593534",
    "2021-10-05 13:39:27.593537",
    "2021-10-05 13:39:27.593540",
    "2021-10-05 13:39:27.593543",
    "2021-10-05 13:39:27.593546",
    "2021-10-05 13:39:27.593549",
    "2021-10-05 13:39:27.593552",
    "2021-10-05 13:39:27.593556",
    "2021-10-05 13:39:27.593559",
    "2021-10-05 13:39:27.593562",
    "2021-10-05 13:39:27.593564",
    "2021-10-05 13:39:27.593567",
    "2021-10-05 13:39:27.593570",
    "2021-10-05 13:39:27.593574",
    "2021-10-05 13:39:27.593577",
    "2021-10-05 13:39:27.593580",
    "2021-10-05 13:39:27.593583",
    "2021-10-05 13:39:27.593586",
    "2021-10-05 13:39:27.593589",
    "2021-10-05 13:39:27.593592",
    "2021-10-05 13:39:27.593595",
    "2021-10-05 13:39:27.593598",
    "2021-10-05 13:39:27.593602",
    "2021-10-05 13:39:27.593605",
    "2021-10-05 13:39:27.593608",
    "2021-10-05 13:39:27.593611",
    "2021-10-05 13:39:27.593614",
    "2021-10-05 13:39:27.593617",
    "2021-10-05 13:39:27.593620",
    "2021-10-05 13:39:27.593623",
    "2021-10-05 13:39:27.593626",
    "2021-10-05 13:39:27.593629",
    "2021-10-05 13:39:27.593631",
    "2021-10-05 13:39:27.593634",
    "2021-10-05 13:39:27.593638",
    "2021-10-05 13:39:27.593641",
    "2021-10-05 13:39:27.593644",
    "2021-10-05 13:39:27.593647",
    "2021-10-05 13:39:27.593650",
    "2021-10-05 13:39:27.593653",
    "2021-10-05 13:39:27.593656",
    "2021-10-05 13:39:27.593659",
    "2021-10-05 13:39:27.593666",
    "2021-10-05 13:39:27.593670",
    "2021-10-05 13:39:27.593673",
    "2021-10-05 13:39:27.593676"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
